# Adding gear-specific harvest outcomes by numbers and weight
#
# - "Eqlb" sheet gets two new data rows: bo / 520000000 and h (reusing the
#   existing "h" shared string already used on the Recruitment sheet).
# - B1 on Eqlb picks up the same scientific-notation number format already
#   used for the big SSB numbers on Recruitment (style index 1 / 0.00E+00).
# - Selection/active-tab bookkeeping flips: Recruitment (not Eqlb) becomes
#   the tab that's active/selected when the workbook is saved, with A13
#   selected there; Eqlb ends up with A3 selected instead of G8.

$wb  = $excel.ActiveWorkbook
$wsRecruitment = $wb.Worksheets.Item("Recruitment")
$wsEqlb        = $wb.Worksheets.Item("Eqlb")

# --- Populate the new data on the Eqlb sheet ---------------------------
$wsEqlb.Activate()

$wsEqlb.Range("A1").Value = "bo"
$wsEqlb.Range("B1").Value = 520000000
$wsEqlb.Range("B1").NumberFormat = $wsRecruitment.Range("B2").NumberFormat

$wsEqlb.Range("A2").Value = "h"

[void]$wsEqlb.Range("A3").Select()

# --- Leave Recruitment as the active sheet/selection --------------------
$wsRecruitment.Activate()
[void]$wsRecruitment.Range("A13").Select()
